# Update "want to go" counts (column F) across sheets, reflecting a fresh
# data scrape (gh-pages output regenerated at commit 456a3b4).
#
# Sheet 1 "展览" (Exhibitions)
# Sheet 3 "本地生活" (Local life)
# Sheet 4 "全部类型" (All types - combined roll-up of the other sheets)

$wb = $excel.ActiveWorkbook

$wsExhibitions = $wb.Worksheets.Item(1)
$wsExhibitions.Range("F6").Value  = 247
$wsExhibitions.Range("F7").Value  = 5021
$wsExhibitions.Range("F9").Value  = 335
$wsExhibitions.Range("F26").Value = 5983
$wsExhibitions.Range("F29").Value = 3212
$wsExhibitions.Range("F30").Value = 321
$wsExhibitions.Range("F34").Value = 115

$wsLocalLife = $wb.Worksheets.Item(3)
$wsLocalLife.Range("F3").Value = 1107

$wsAllTypes = $wb.Worksheets.Item(4)
$wsAllTypes.Range("F4").Value  = 1107
$wsAllTypes.Range("F10").Value = 247
$wsAllTypes.Range("F11").Value = 5021
$wsAllTypes.Range("F13").Value = 335
$wsAllTypes.Range("F30").Value = 5983
$wsAllTypes.Range("F33").Value = 3212
$wsAllTypes.Range("F34").Value = 321
$wsAllTypes.Range("F39").Value = 115
